$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (cell style) from the last existing data row down through the new rows
$ws.Range("A587:C587").Copy() | Out-Null
$ws.Range("A588:C736").PasteSpecial(-4122) | Out-Null
$ws.Range("A588:C736").RowHeight = 27.6

# Populate the new registration/course/grade rows
$ws.Cells.Item(588, 1).Value = "ME231006"
$ws.Cells.Item(588, 2).Value = "ME-2108L"
$ws.Cells.Item(588, 3).Value = "C+"
$ws.Cells.Item(589, 1).Value = "ME231006"
$ws.Cells.Item(589, 2).Value = "ME-2303"
$ws.Cells.Item(589, 3).Value = "B-"
$ws.Cells.Item(590, 1).Value = "ME231006"
$ws.Cells.Item(590, 2).Value = "ME-2305L"
$ws.Cells.Item(590, 3).Value = "B+"
$ws.Cells.Item(591, 1).Value = "ME231006"
$ws.Cells.Item(591, 2).Value = "ME-2105"
$ws.Cells.Item(591, 3).Value = "B"
$ws.Cells.Item(592, 1).Value = "ME231006"
$ws.Cells.Item(592, 2).Value = "ME-2402"
$ws.Cells.Item(592, 3).Value = "D"
$ws.Cells.Item(593, 1).Value = "ME231006"
$ws.Cells.Item(593, 2).Value = "ME-2402L"
$ws.Cells.Item(593, 3).Value = "B"
$ws.Cells.Item(594, 1).Value = "ME231006"
$ws.Cells.Item(594, 2).Value = "ME-2103"
$ws.Cells.Item(594, 3).Value = "D+"
$ws.Cells.Item(595, 1).Value = "ME231006"
$ws.Cells.Item(595, 2).Value = "ME-2107L"
$ws.Cells.Item(595, 3).Value = "A-"
$ws.Cells.Item(596, 1).Value = "ME231006"
$ws.Cells.Item(596, 2).Value = "BS-1402"
$ws.Cells.Item(596, 3).Value = "C-"
$ws.Cells.Item(597, 1).Value = "ME231007"
$ws.Cells.Item(597, 2).Value = "ME-2108L"
$ws.Cells.Item(597, 3).Value = "B+"
$ws.Cells.Item(598, 1).Value = "ME231007"
$ws.Cells.Item(598, 2).Value = "ME-2303"
$ws.Cells.Item(598, 3).Value = "C"
$ws.Cells.Item(599, 1).Value = "ME231007"
$ws.Cells.Item(599, 2).Value = "ME-2305L"
$ws.Cells.Item(599, 3).Value = "B"
$ws.Cells.Item(600, 1).Value = "ME231007"
$ws.Cells.Item(600, 2).Value = "ME-2105"
$ws.Cells.Item(600, 3).Value = "C"
$ws.Cells.Item(601, 1).Value = "ME231007"
$ws.Cells.Item(601, 2).Value = "ME-2402"
$ws.Cells.Item(601, 3).Value = "C+"
$ws.Cells.Item(602, 1).Value = "ME231007"
$ws.Cells.Item(602, 2).Value = "ME-2402L"
$ws.Cells.Item(602, 3).Value = "B+"
$ws.Cells.Item(603, 1).Value = "ME231007"
$ws.Cells.Item(603, 2).Value = "ME-2103"
$ws.Cells.Item(603, 3).Value = "F"
$ws.Cells.Item(604, 1).Value = "ME231007"
$ws.Cells.Item(604, 2).Value = "ME-2107L"
$ws.Cells.Item(604, 3).Value = "B+"
$ws.Cells.Item(605, 1).Value = "ME231007"
$ws.Cells.Item(605, 2).Value = "BS-1402"
$ws.Cells.Item(605, 3).Value = "D"
$ws.Cells.Item(606, 1).Value = "ME231009"
$ws.Cells.Item(606, 2).Value = "ME-2108L"
$ws.Cells.Item(606, 3).Value = "B+"
$ws.Cells.Item(607, 1).Value = "ME231009"
$ws.Cells.Item(607, 2).Value = "ME-2303"
$ws.Cells.Item(607, 3).Value = "A"
$ws.Cells.Item(608, 1).Value = "ME231009"
$ws.Cells.Item(608, 2).Value = "ME-2305L"
$ws.Cells.Item(608, 3).Value = "B-"
$ws.Cells.Item(609, 1).Value = "ME231009"
$ws.Cells.Item(609, 2).Value = "ME-2105"
$ws.Cells.Item(609, 3).Value = "C"
$ws.Cells.Item(610, 1).Value = "ME231009"
$ws.Cells.Item(610, 2).Value = "ME-2402"
$ws.Cells.Item(610, 3).Value = "B-"
$ws.Cells.Item(611, 1).Value = "ME231009"
$ws.Cells.Item(611, 2).Value = "ME-2402L"
$ws.Cells.Item(611, 3).Value = "A-"
$ws.Cells.Item(612, 1).Value = "ME231009"
$ws.Cells.Item(612, 2).Value = "ME-2103"
$ws.Cells.Item(612, 3).Value = "C"
$ws.Cells.Item(613, 1).Value = "ME231009"
$ws.Cells.Item(613, 2).Value = "ME-2107L"
$ws.Cells.Item(613, 3).Value = "A-"
$ws.Cells.Item(614, 1).Value = "ME231009"
$ws.Cells.Item(614, 2).Value = "BS-1402"
$ws.Cells.Item(614, 3).Value = "A"
$ws.Cells.Item(615, 1).Value = "ME231010"
$ws.Cells.Item(615, 2).Value = "ME-2108L"
$ws.Cells.Item(615, 3).Value = "A-"
$ws.Cells.Item(616, 1).Value = "ME231010"
$ws.Cells.Item(616, 2).Value = "ME-2303"
$ws.Cells.Item(616, 3).Value = "F"
$ws.Cells.Item(617, 1).Value = "ME231010"
$ws.Cells.Item(617, 2).Value = "ME-2305L"
$ws.Cells.Item(617, 3).Value = "A"
$ws.Cells.Item(618, 1).Value = "ME231010"
$ws.Cells.Item(618, 2).Value = "ME-2105"
$ws.Cells.Item(618, 3).Value = "D"
$ws.Cells.Item(619, 1).Value = "ME231010"
$ws.Cells.Item(619, 2).Value = "ME-2402"
$ws.Cells.Item(619, 3).Value = "C+"
$ws.Cells.Item(620, 1).Value = "ME231010"
$ws.Cells.Item(620, 2).Value = "ME-2402L"
$ws.Cells.Item(620, 3).Value = "B+"
$ws.Cells.Item(621, 1).Value = "ME231010"
$ws.Cells.Item(621, 2).Value = "ME-2103"
$ws.Cells.Item(621, 3).Value = "D"
$ws.Cells.Item(622, 1).Value = "ME231010"
$ws.Cells.Item(622, 2).Value = "ME-2107L"
$ws.Cells.Item(622, 3).Value = "A-"
$ws.Cells.Item(623, 1).Value = "ME231010"
$ws.Cells.Item(623, 2).Value = "BS-1402"
$ws.Cells.Item(623, 3).Value = "C"
$ws.Cells.Item(624, 1).Value = "ME231011"
$ws.Cells.Item(624, 2).Value = "ME-2108L"
$ws.Cells.Item(624, 3).Value = "C-"
$ws.Cells.Item(625, 1).Value = "ME231011"
$ws.Cells.Item(625, 2).Value = "ME-1103"
$ws.Cells.Item(625, 3).Value = "C+"
$ws.Cells.Item(626, 1).Value = "ME231011"
$ws.Cells.Item(626, 2).Value = "ME-2402"
$ws.Cells.Item(626, 3).Value = "F"
$ws.Cells.Item(627, 1).Value = "ME231011"
$ws.Cells.Item(627, 2).Value = "ME-2402L"
$ws.Cells.Item(627, 3).Value = "B"
$ws.Cells.Item(628, 1).Value = "ME231011"
$ws.Cells.Item(628, 2).Value = "BS-1402"
$ws.Cells.Item(628, 3).Value = "F"
$ws.Cells.Item(629, 1).Value = "ME231012"
$ws.Cells.Item(629, 2).Value = "ME-2108L"
$ws.Cells.Item(629, 3).Value = "B+"
$ws.Cells.Item(630, 1).Value = "ME231012"
$ws.Cells.Item(630, 2).Value = "ME-2303"
$ws.Cells.Item(630, 3).Value = "B-"
$ws.Cells.Item(631, 1).Value = "ME231012"
$ws.Cells.Item(631, 2).Value = "ME-2305L"
$ws.Cells.Item(631, 3).Value = "A"
$ws.Cells.Item(632, 1).Value = "ME231012"
$ws.Cells.Item(632, 2).Value = "ME-2105"
$ws.Cells.Item(632, 3).Value = "B-"
$ws.Cells.Item(633, 1).Value = "ME231012"
$ws.Cells.Item(633, 2).Value = "ME-2402"
$ws.Cells.Item(633, 3).Value = "C+"
$ws.Cells.Item(634, 1).Value = "ME231012"
$ws.Cells.Item(634, 2).Value = "ME-2402L"
$ws.Cells.Item(634, 3).Value = "B+"
$ws.Cells.Item(635, 1).Value = "ME231012"
$ws.Cells.Item(635, 2).Value = "ME-2103"
$ws.Cells.Item(635, 3).Value = "C"
$ws.Cells.Item(636, 1).Value = "ME231012"
$ws.Cells.Item(636, 2).Value = "ME-2107L"
$ws.Cells.Item(636, 3).Value = "B+"
$ws.Cells.Item(637, 1).Value = "ME231012"
$ws.Cells.Item(637, 2).Value = "BS-1402"
$ws.Cells.Item(637, 3).Value = "B"
$ws.Cells.Item(638, 1).Value = "ME231013"
$ws.Cells.Item(638, 2).Value = "ME-2108L"
$ws.Cells.Item(638, 3).Value = "B-"
$ws.Cells.Item(639, 1).Value = "ME231013"
$ws.Cells.Item(639, 2).Value = "ME-2303"
$ws.Cells.Item(639, 3).Value = "D"
$ws.Cells.Item(640, 1).Value = "ME231013"
$ws.Cells.Item(640, 2).Value = "ME-2305L"
$ws.Cells.Item(640, 3).Value = "B"
$ws.Cells.Item(641, 1).Value = "ME231013"
$ws.Cells.Item(641, 2).Value = "BS-2303"
$ws.Cells.Item(641, 3).Value = "F"
$ws.Cells.Item(642, 1).Value = "ME231013"
$ws.Cells.Item(642, 2).Value = "ME-2105"
$ws.Cells.Item(642, 3).Value = "D"
$ws.Cells.Item(643, 1).Value = "ME231013"
$ws.Cells.Item(643, 2).Value = "ME-2402"
$ws.Cells.Item(643, 3).Value = "C+"
$ws.Cells.Item(644, 1).Value = "ME231013"
$ws.Cells.Item(644, 2).Value = "ME-2402L"
$ws.Cells.Item(644, 3).Value = "B"
$ws.Cells.Item(645, 1).Value = "ME231013"
$ws.Cells.Item(645, 2).Value = "ME-2103"
$ws.Cells.Item(645, 3).Value = "D"
$ws.Cells.Item(646, 1).Value = "ME231013"
$ws.Cells.Item(646, 2).Value = "ME-2107L"
$ws.Cells.Item(646, 3).Value = "B+"
$ws.Cells.Item(647, 1).Value = "ME231015"
$ws.Cells.Item(647, 2).Value = "ME-2108L"
$ws.Cells.Item(647, 3).Value = "A-"
$ws.Cells.Item(648, 1).Value = "ME231015"
$ws.Cells.Item(648, 2).Value = "ME-2303"
$ws.Cells.Item(648, 3).Value = "A-"
$ws.Cells.Item(649, 1).Value = "ME231015"
$ws.Cells.Item(649, 2).Value = "ME-2305L"
$ws.Cells.Item(649, 3).Value = "B"
$ws.Cells.Item(650, 1).Value = "ME231015"
$ws.Cells.Item(650, 2).Value = "ME-2105"
$ws.Cells.Item(650, 3).Value = "B+"
$ws.Cells.Item(651, 1).Value = "ME231015"
$ws.Cells.Item(651, 2).Value = "ME-2402"
$ws.Cells.Item(651, 3).Value = "B"
$ws.Cells.Item(652, 1).Value = "ME231015"
$ws.Cells.Item(652, 2).Value = "ME-2402L"
$ws.Cells.Item(652, 3).Value = "A"
$ws.Cells.Item(653, 1).Value = "ME231015"
$ws.Cells.Item(653, 2).Value = "ME-2103"
$ws.Cells.Item(653, 3).Value = "B"
$ws.Cells.Item(654, 1).Value = "ME231015"
$ws.Cells.Item(654, 2).Value = "ME-2107L"
$ws.Cells.Item(654, 3).Value = "A"
$ws.Cells.Item(655, 1).Value = "ME231015"
$ws.Cells.Item(655, 2).Value = "BS-1402"
$ws.Cells.Item(655, 3).Value = "A"
$ws.Cells.Item(656, 1).Value = "ME231016"
$ws.Cells.Item(656, 2).Value = "ME-2108L"
$ws.Cells.Item(656, 3).Value = "C"
$ws.Cells.Item(657, 1).Value = "ME231016"
$ws.Cells.Item(657, 2).Value = "ME-1103"
$ws.Cells.Item(657, 3).Value = "F"
$ws.Cells.Item(658, 1).Value = "ME231016"
$ws.Cells.Item(658, 2).Value = "ME-2303"
$ws.Cells.Item(658, 3).Value = "F"
$ws.Cells.Item(659, 1).Value = "ME231016"
$ws.Cells.Item(659, 2).Value = "ME-2305L"
$ws.Cells.Item(659, 3).Value = "A"
$ws.Cells.Item(660, 1).Value = "ME231016"
$ws.Cells.Item(660, 2).Value = "BS-2303"
$ws.Cells.Item(660, 3).Value = "F"
$ws.Cells.Item(661, 1).Value = "ME231016"
$ws.Cells.Item(661, 2).Value = "ME-2402"
$ws.Cells.Item(661, 3).Value = "F"
$ws.Cells.Item(662, 1).Value = "ME231016"
$ws.Cells.Item(662, 2).Value = "ME-2402L"
$ws.Cells.Item(662, 3).Value = "C+"
$ws.Cells.Item(663, 1).Value = "ME231016"
$ws.Cells.Item(663, 2).Value = "ME-1301"
$ws.Cells.Item(663, 3).Value = "F"
$ws.Cells.Item(664, 1).Value = "ME231017"
$ws.Cells.Item(664, 2).Value = "ME-2108L"
$ws.Cells.Item(664, 3).Value = "D+"
$ws.Cells.Item(665, 1).Value = "ME231017"
$ws.Cells.Item(665, 2).Value = "ME-1103"
$ws.Cells.Item(665, 3).Value = "F"
$ws.Cells.Item(666, 1).Value = "ME231017"
$ws.Cells.Item(666, 2).Value = "ME-2303"
$ws.Cells.Item(666, 3).Value = "D"
$ws.Cells.Item(667, 1).Value = "ME231017"
$ws.Cells.Item(667, 2).Value = "ME-2305L"
$ws.Cells.Item(667, 3).Value = "A"
$ws.Cells.Item(668, 1).Value = "ME231017"
$ws.Cells.Item(668, 2).Value = "ME-2402"
$ws.Cells.Item(668, 3).Value = "C"
$ws.Cells.Item(669, 1).Value = "ME231017"
$ws.Cells.Item(669, 2).Value = "ME-2402L"
$ws.Cells.Item(669, 3).Value = "C+"
$ws.Cells.Item(670, 1).Value = "ME231017"
$ws.Cells.Item(670, 2).Value = "BS-1402"
$ws.Cells.Item(670, 3).Value = "D"
$ws.Cells.Item(671, 1).Value = "ME231018"
$ws.Cells.Item(671, 2).Value = "ME-2108L"
$ws.Cells.Item(671, 3).Value = "A-"
$ws.Cells.Item(672, 1).Value = "ME231018"
$ws.Cells.Item(672, 2).Value = "ME-2303"
$ws.Cells.Item(672, 3).Value = "C+"
$ws.Cells.Item(673, 1).Value = "ME231018"
$ws.Cells.Item(673, 2).Value = "ME-2305L"
$ws.Cells.Item(673, 3).Value = "A"
$ws.Cells.Item(674, 1).Value = "ME231018"
$ws.Cells.Item(674, 2).Value = "ME-2105"
$ws.Cells.Item(674, 3).Value = "D"
$ws.Cells.Item(675, 1).Value = "ME231018"
$ws.Cells.Item(675, 2).Value = "ME-2402"
$ws.Cells.Item(675, 3).Value = "C"
$ws.Cells.Item(676, 1).Value = "ME231018"
$ws.Cells.Item(676, 2).Value = "ME-2402L"
$ws.Cells.Item(676, 3).Value = "A-"
$ws.Cells.Item(677, 1).Value = "ME231018"
$ws.Cells.Item(677, 2).Value = "ME-2103"
$ws.Cells.Item(677, 3).Value = "F"
$ws.Cells.Item(678, 1).Value = "ME231018"
$ws.Cells.Item(678, 2).Value = "ME-2107L"
$ws.Cells.Item(678, 3).Value = "B-"
$ws.Cells.Item(679, 1).Value = "ME231018"
$ws.Cells.Item(679, 2).Value = "BS-1402"
$ws.Cells.Item(679, 3).Value = "D"
$ws.Cells.Item(680, 1).Value = "ME231022"
$ws.Cells.Item(680, 2).Value = "ME-2108L"
$ws.Cells.Item(680, 3).Value = "A-"
$ws.Cells.Item(681, 1).Value = "ME231022"
$ws.Cells.Item(681, 2).Value = "ME-2303"
$ws.Cells.Item(681, 3).Value = "B-"
$ws.Cells.Item(682, 1).Value = "ME231022"
$ws.Cells.Item(682, 2).Value = "ME-2305L"
$ws.Cells.Item(682, 3).Value = "A"
$ws.Cells.Item(683, 1).Value = "ME231022"
$ws.Cells.Item(683, 2).Value = "ME-2105"
$ws.Cells.Item(683, 3).Value = "C-"
$ws.Cells.Item(684, 1).Value = "ME231022"
$ws.Cells.Item(684, 2).Value = "ME-2402"
$ws.Cells.Item(684, 3).Value = "B"
$ws.Cells.Item(685, 1).Value = "ME231022"
$ws.Cells.Item(685, 2).Value = "ME-2402L"
$ws.Cells.Item(685, 3).Value = "B-"
$ws.Cells.Item(686, 1).Value = "ME231022"
$ws.Cells.Item(686, 2).Value = "ME-2103"
$ws.Cells.Item(686, 3).Value = "D"
$ws.Cells.Item(687, 1).Value = "ME231022"
$ws.Cells.Item(687, 2).Value = "ME-2107L"
$ws.Cells.Item(687, 3).Value = "A-"
$ws.Cells.Item(688, 1).Value = "ME231022"
$ws.Cells.Item(688, 2).Value = "BS-1402"
$ws.Cells.Item(688, 3).Value = "C"
$ws.Cells.Item(689, 1).Value = "ME231023"
$ws.Cells.Item(689, 2).Value = "ME-2108L"
$ws.Cells.Item(689, 3).Value = "B"
$ws.Cells.Item(690, 1).Value = "ME231023"
$ws.Cells.Item(690, 2).Value = "CS-1502"
$ws.Cells.Item(690, 3).Value = "F"
$ws.Cells.Item(691, 1).Value = "ME231023"
$ws.Cells.Item(691, 2).Value = "ME-1103"
$ws.Cells.Item(691, 3).Value = "F"
$ws.Cells.Item(692, 1).Value = "ME231023"
$ws.Cells.Item(692, 2).Value = "ME-2303"
$ws.Cells.Item(692, 3).Value = "F"
$ws.Cells.Item(693, 1).Value = "ME231023"
$ws.Cells.Item(693, 2).Value = "ME-2305L"
$ws.Cells.Item(693, 3).Value = "C+"
$ws.Cells.Item(694, 1).Value = "ME231023"
$ws.Cells.Item(694, 2).Value = "ME-2402"
$ws.Cells.Item(694, 3).Value = "F"
$ws.Cells.Item(695, 1).Value = "ME231023"
$ws.Cells.Item(695, 2).Value = "ME-2402L"
$ws.Cells.Item(695, 3).Value = "B"
$ws.Cells.Item(696, 1).Value = "ME231026"
$ws.Cells.Item(696, 2).Value = "ME-2108L"
$ws.Cells.Item(696, 3).Value = "D+"
$ws.Cells.Item(697, 1).Value = "ME231026"
$ws.Cells.Item(697, 2).Value = "ME-1103"
$ws.Cells.Item(697, 3).Value = "F"
$ws.Cells.Item(698, 1).Value = "ME231026"
$ws.Cells.Item(698, 2).Value = "ME-2402"
$ws.Cells.Item(698, 3).Value = "D"
$ws.Cells.Item(699, 1).Value = "ME231026"
$ws.Cells.Item(699, 2).Value = "ME-2402L"
$ws.Cells.Item(699, 3).Value = "C"
$ws.Cells.Item(700, 1).Value = "ME231026"
$ws.Cells.Item(700, 2).Value = "BS-1402"
$ws.Cells.Item(700, 3).Value = "D+"
$ws.Cells.Item(701, 1).Value = "ME231031"
$ws.Cells.Item(701, 2).Value = "ME-2108L"
$ws.Cells.Item(701, 3).Value = "A-"
$ws.Cells.Item(702, 1).Value = "ME231031"
$ws.Cells.Item(702, 2).Value = "ME-2303"
$ws.Cells.Item(702, 3).Value = "A"
$ws.Cells.Item(703, 1).Value = "ME231031"
$ws.Cells.Item(703, 2).Value = "ME-2305L"
$ws.Cells.Item(703, 3).Value = "A"
$ws.Cells.Item(704, 1).Value = "ME231031"
$ws.Cells.Item(704, 2).Value = "ME-2105"
$ws.Cells.Item(704, 3).Value = "B+"
$ws.Cells.Item(705, 1).Value = "ME231031"
$ws.Cells.Item(705, 2).Value = "ME-2402"
$ws.Cells.Item(705, 3).Value = "B-"
$ws.Cells.Item(706, 1).Value = "ME231031"
$ws.Cells.Item(706, 2).Value = "ME-2402L"
$ws.Cells.Item(706, 3).Value = "B+"
$ws.Cells.Item(707, 1).Value = "ME231031"
$ws.Cells.Item(707, 2).Value = "ME-2103"
$ws.Cells.Item(707, 3).Value = "C"
$ws.Cells.Item(708, 1).Value = "ME231031"
$ws.Cells.Item(708, 2).Value = "ME-2107L"
$ws.Cells.Item(708, 3).Value = "A-"
$ws.Cells.Item(709, 1).Value = "ME231031"
$ws.Cells.Item(709, 2).Value = "BS-1402"
$ws.Cells.Item(709, 3).Value = "A-"
$ws.Cells.Item(710, 1).Value = "ME231035"
$ws.Cells.Item(710, 2).Value = "ME-2108L"
$ws.Cells.Item(710, 3).Value = "A-"
$ws.Cells.Item(711, 1).Value = "ME231035"
$ws.Cells.Item(711, 2).Value = "ME-2303"
$ws.Cells.Item(711, 3).Value = "A-"
$ws.Cells.Item(712, 1).Value = "ME231035"
$ws.Cells.Item(712, 2).Value = "ME-2305L"
$ws.Cells.Item(712, 3).Value = "B+"
$ws.Cells.Item(713, 1).Value = "ME231035"
$ws.Cells.Item(713, 2).Value = "ME-2105"
$ws.Cells.Item(713, 3).Value = "B+"
$ws.Cells.Item(714, 1).Value = "ME231035"
$ws.Cells.Item(714, 2).Value = "ME-2402"
$ws.Cells.Item(714, 3).Value = "B+"
$ws.Cells.Item(715, 1).Value = "ME231035"
$ws.Cells.Item(715, 2).Value = "ME-2402L"
$ws.Cells.Item(715, 3).Value = "A-"
$ws.Cells.Item(716, 1).Value = "ME231035"
$ws.Cells.Item(716, 2).Value = "ME-2103"
$ws.Cells.Item(716, 3).Value = "C"
$ws.Cells.Item(717, 1).Value = "ME231035"
$ws.Cells.Item(717, 2).Value = "ME-2107L"
$ws.Cells.Item(717, 3).Value = "A"
$ws.Cells.Item(718, 1).Value = "ME231035"
$ws.Cells.Item(718, 2).Value = "BS-1402"
$ws.Cells.Item(718, 3).Value = "A-"
$ws.Cells.Item(719, 1).Value = "ME231037"
$ws.Cells.Item(719, 2).Value = "ME-2108L"
$ws.Cells.Item(719, 3).Value = "B+"
$ws.Cells.Item(720, 1).Value = "ME231037"
$ws.Cells.Item(720, 2).Value = "ME-2303"
$ws.Cells.Item(720, 3).Value = "D+"
$ws.Cells.Item(721, 1).Value = "ME231037"
$ws.Cells.Item(721, 2).Value = "ME-2305L"
$ws.Cells.Item(721, 3).Value = "A-"
$ws.Cells.Item(722, 1).Value = "ME231037"
$ws.Cells.Item(722, 2).Value = "BS-2303"
$ws.Cells.Item(722, 3).Value = "B-"
$ws.Cells.Item(723, 1).Value = "ME231037"
$ws.Cells.Item(723, 2).Value = "ME-2105"
$ws.Cells.Item(723, 3).Value = "C"
$ws.Cells.Item(724, 1).Value = "ME231037"
$ws.Cells.Item(724, 2).Value = "ME-2402"
$ws.Cells.Item(724, 3).Value = "C"
$ws.Cells.Item(725, 1).Value = "ME231037"
$ws.Cells.Item(725, 2).Value = "ME-2402L"
$ws.Cells.Item(725, 3).Value = "B-"
$ws.Cells.Item(726, 1).Value = "ME231037"
$ws.Cells.Item(726, 2).Value = "ME-2103"
$ws.Cells.Item(726, 3).Value = "F"
$ws.Cells.Item(727, 1).Value = "ME231037"
$ws.Cells.Item(727, 2).Value = "ME-2107L"
$ws.Cells.Item(727, 3).Value = "C+"
$ws.Cells.Item(728, 1).Value = "ME231038"
$ws.Cells.Item(728, 2).Value = "ME-2108L"
$ws.Cells.Item(728, 3).Value = "A"
$ws.Cells.Item(729, 1).Value = "ME231038"
$ws.Cells.Item(729, 2).Value = "ME-2303"
$ws.Cells.Item(729, 3).Value = "A-"
$ws.Cells.Item(730, 1).Value = "ME231038"
$ws.Cells.Item(730, 2).Value = "ME-2305L"
$ws.Cells.Item(730, 3).Value = "B+"
$ws.Cells.Item(731, 1).Value = "ME231038"
$ws.Cells.Item(731, 2).Value = "ME-2105"
$ws.Cells.Item(731, 3).Value = "B-"
$ws.Cells.Item(732, 1).Value = "ME231038"
$ws.Cells.Item(732, 2).Value = "ME-2402"
$ws.Cells.Item(732, 3).Value = "B-"
$ws.Cells.Item(733, 1).Value = "ME231038"
$ws.Cells.Item(733, 2).Value = "ME-2402L"
$ws.Cells.Item(733, 3).Value = "A"
$ws.Cells.Item(734, 1).Value = "ME231038"
$ws.Cells.Item(734, 2).Value = "ME-2103"
$ws.Cells.Item(734, 3).Value = "C-"
$ws.Cells.Item(735, 1).Value = "ME231038"
$ws.Cells.Item(735, 2).Value = "ME-2107L"
$ws.Cells.Item(735, 3).Value = "A"
$ws.Cells.Item(736, 1).Value = "ME231038"
$ws.Cells.Item(736, 2).Value = "BS-1402"
$ws.Cells.Item(736, 3).Value = "C-"

# Match the final selection left by the author
$ws.Range("A588:C736").Select() | Out-Null

Write-Host "Added 149 rows (588-736)"
